$d = $word.ActiveDocument

# Remove the first paragraph entirely (the "I personally examined..." note),
# including its paragraph mark, so the document now starts with "OBJECTIVE:".
$p = $d.Paragraphs.Item(1)
$p.Range.Delete()
